# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the newly generated data snapshot from the gh-pages build.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 3703
$ws1.Range("F10").Value = 2961
$ws1.Range("F13").Value = 2257
$ws1.Range("F16").Value = 39
$ws1.Range("F19").Value = 184
$ws1.Range("F28").Value = 141
$ws1.Range("F30").Value = 4086
$ws1.Range("F31").Value = 3651
$ws1.Range("F32").Value = 55
$ws1.Range("F34").Value = 1088
$ws1.Range("F37").Value = 1300

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2194

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 3703
$ws4.Range("F15").Value = 2961
$ws4.Range("F17").Value = 2257
$ws4.Range("F20").Value = 39
$ws4.Range("F23").Value = 184
$ws4.Range("F30").Value = 141
$ws4.Range("F33").Value = 4086
$ws4.Range("F34").Value = 3651
$ws4.Range("F35").Value = 55
$ws4.Range("F36").Value = 1088
$ws4.Range("F43").Value = 1300
